$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 13:29"

# Re-sorted province/city table (descending by Casos totales) with updated counts
# for Zaragoza, A Coruna, Pontevedra, Caceres, Badajoz, Ourense, Lugo, Teruel, Huesca
$data = @(
    @("Madrid", 21520, 6326, 12437, 2757),
    @("Cataluña", 14263, 3106, 10087, 1070),
    @("Castilla-La Mancha", 2780, 71, 2446, 263),
    @("Valencia/Valencia", 2263, 68, 2091, 104),
    @("Bizkaia/Vizcaya", 2186, 814, 1793, 79),
    @("Navarra", 1829, 98, 1661, 70),
    @("Araba/Alava", 1610, 814, 1207, 103),
    @("La Rioja", 1436, 364, 1007, 65),
    @("Ciudad Real", 1422, 153, 1273, 89),
    @("Alacant/Alicante", 1322, 20, 1193, 109),
    @("Zaragoza", 1239, 96, 1071, 72),
    @("A Coruña", 1177, 95, 1114, 31),
    @("Toledo", 1112, 153, 972, 90),
    @("Gran Canaria", 1025, 20, 834, 36),
    @("Malaga", 1006, 61, 899, 46),
    @("Asturias", 1004, 65, 906, 33),
    @("Pontevedra", 960, 95, 923, 9),
    @("Caceres", 957, 38, 841, 78),
    @("Cantabria", 937, 21, 894, 22),
    @("Aragon", 907, 29, 838, 40),
    @("Gipuzkoa/Guipuzcoa", 805, 814, 580, 25),
    @("Murcia", 802, 12, 773, 17),
    @("Albacete", 780, 153, 667, 83),
    @("Granada", 711, 1, 676, 34),
    @("Sevilla", 708, 8, 675, 25),
    @("Salamanca", 629, 73, 502, 54),
    @("Valladolid", 598, 50, 521, 27),
    @("Tenerife", 539, 15, 519, 36),
    @("Burgos", 530, 67, 431, 32),
    @("Leon", 516, 48, 422, 46),
    @("Castello/Castellon", 449, 4, 424, 21),
    @("Guadalajara", 440, 153, 362, 75),
    @("Badajoz", 437, 11, 416, 10),
    @("Jaen", 414, 7, 392, 15),
    @("Ourense", 396, 95, 371, 6),
    @("Segovia", 361, 62, 262, 37),
    @("Cordoba", 359, 1, 350, 8),
    @("Cadiz", 348, 6, 336, 6),
    @("Soria", 339, 32, 291, 16),
    @("Avila", 270, 53, 190, 27),
    @("Lugo", 239, 95, 219, 4),
    @("Mallorca", 210, 18, 194, 12),
    @("Cuenca", 180, 153, 130, 40),
    @("Teruel", 179, 10, 159, 10),
    @("Huesca", 174, 14, 151, 9),
    @("Almeria", 142, 6, 128, 8),
    @("Palencia", 139, 14, 120, 5),
    @("Zamora", 106, 24, 73, 9),
    @("Huelva", 105, 2, 101, 2),
    @("Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3),
    @("Melilla", 45, 0, 44, 1),
    @("La Palma", 33, 15, 33, 36),
    @("Ibiza", 21, 18, 20, 1),
    @("Fuerteventura", 20, 15, 20, 36),
    @("Lanzarote", 17, 15, 17, 36),
    @("Ceuta", 17, 0, 16, 1),
    @("Menorca", 15, 18, 13, 0),
    @("Arroyo de la Luz", 7, 0, 7, 0),
    @("La Gomera", 4, 15, 2, 36),
    @("El Hierro", 3, 15, 3, 36),
    @("Formentera", 0, 10, 0, 8)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
